$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: preserve the previous "Ultimate Load" test-column values that
# used to live in column B.
$ws.Range("C1").Value = "Ultimate Load"
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 4680
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 120
$ws.Range("C7").Value = 60

# Column B: new "Checkout Load" test data.
$ws.Range("B1").Value = "Checkout Load"
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 800
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 30
$ws.Range("B7").Value = 10

$ws.Range("B4").Select() | Out-Null
